$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

$ws.Range("A4").Value = "에이치브이엠(구.한국진공야금)"
$ws.Range("B4").Value = "2024.06.11~06.17"
$ws.Range("C4").Value = "11,000~14,200"
$ws.Range("D4").Value = "-"
$ws.Range("E4").Value = "26400"
$ws.Range("F4").Value = "NH투자증권"

$ws.Range("A5").Value = "이노스페이스"
$ws.Range("B5").Value = "2024.06.11~06.17"
$ws.Range("C5").Value = "36,400~43,300"
$ws.Range("D5").Value = "-"
$ws.Range("E5").Value = "48412"
$ws.Range("F5").Value = "미래에셋증권,신한투자증권"

$ws.Range("A6").Value = "한국스팩15호"
$ws.Range("B6").Value = "2024.06.10~06.11"
$ws.Range("C6").Value = "2,000~2,000"
$ws.Range("D6").Value = "-"
$ws.Range("E6").Value = "12500"
$ws.Range("F6").Value = "한국투자증권"

$ws.Range("A7").Value = "하이젠알앤엠"
$ws.Range("B7").Value = "2024.06.07~06.13"
$ws.Range("C7").Value = "4,500~5,500"
$ws.Range("D7").Value = "-"
$ws.Range("E7").Value = "15300"
$ws.Range("F7").Value = "한국투자증권"

$ws.Range("A8").Value = "미래에셋비전스팩6호"
$ws.Range("B8").Value = "2024.06.05~06.07"
$ws.Range("C8").Value = "2,000~2,000"
$ws.Range("D8").Value = "-"
$ws.Range("E8").Value = "12900"
$ws.Range("F8").Value = "미래에셋증권"

$ws.Range("A9").Value = "KB스팩29호"
$ws.Range("B9").Value = "2024.06.04~06.05"
$ws.Range("C9").Value = "2,000~2,000"
$ws.Range("D9").Value = "-"
$ws.Range("E9").Value = "12000"
$ws.Range("F9").Value = "KB증권"

$ws.Range("A10").Value = "에이치엠씨아이비스팩7호"
$ws.Range("B10").Value = "2024.06.04~06.05"
$ws.Range("C10").Value = "2,000~2,000"
$ws.Range("D10").Value = "-"
$ws.Range("E10").Value = "14000"
$ws.Range("F10").Value = "현대차증권"

$ws.Range("A11").Value = "에스오에스랩"
$ws.Range("B11").Value = "2024.06.03~06.10"
$ws.Range("C11").Value = "7,500~9,000"
$ws.Range("D11").Value = "-"
$ws.Range("E11").Value = "15000"
$ws.Range("F11").Value = "한국투자증권"

$ws.Range("A12").Value = "미래에셋비전스팩5호"
$ws.Range("B12").Value = "2024.06.03~06.04"
$ws.Range("C12").Value = "2,000~2,000"
$ws.Range("D12").Value = "-"
$ws.Range("E12").Value = "9500"
$ws.Range("F12").Value = "미래에셋증권"

$ws.Range("A13").Value = "한국스팩14호"
$ws.Range("B13").Value = "2024.06.03~06.04"
$ws.Range("C13").Value = "2,000~2,000"
$ws.Range("D13").Value = "-"
$ws.Range("E13").Value = "8000"
$ws.Range("F13").Value = "한국투자증권"

$ws.Range("A14").Value = "엑셀세라퓨틱스"
$ws.Range("B14").Value = "2024.06.03~06.10"
$ws.Range("C14").Value = "6,200~7,700"
$ws.Range("D14").Value = "-"
$ws.Range("E14").Value = "10032"
$ws.Range("F14").Value = "대신증권"

$ws.Range("A15").Value = "시프트업"
$ws.Range("B15").Value = "2024.06.03~06.13"
$ws.Range("C15").Value = "47,000~60,000"
$ws.Range("D15").Value = "-"
$ws.Range("E15").Value = "340750"
$ws.Range("F15").Value = "한국투자증권,NH투자증권,신한투자증권"

$ws.Range("A16").Value = "이노그리드"
$ws.Range("B16").Value = "2024.05.31~06.07"
$ws.Range("C16").Value = "29,000~35,000"
$ws.Range("D16").Value = "-"
$ws.Range("E16").Value = "17400"
$ws.Range("F16").Value = "한국투자증권"

$ws.Range("A17").Value = "디비금융스팩12호"
$ws.Range("B17").Value = "2024.05.28~05.29"
$ws.Range("C17").Value = "2,000~2,000"
$ws.Range("D17").Value = "-"
$ws.Range("E17").Value = "10000"
$ws.Range("F17").Value = "DB금융투자"

$ws.Range("A18").Value = "씨어스테크놀로지"
$ws.Range("B18").Value = "2024.05.27~05.31"
$ws.Range("C18").Value = "10,500~14,000"
$ws.Range("D18").Value = "-"
$ws.Range("E18").Value = "13650"
$ws.Range("F18").Value = "한국투자증권"

$ws.Range("A19").Value = "라메디텍"
$ws.Range("B19").Value = "2024.05.27~05.31"
$ws.Range("C19").Value = "10,400~12,700"
$ws.Range("D19").Value = "-"
$ws.Range("E19").Value = "13499"
$ws.Range("F19").Value = "대신증권"

$ws.Range("A20").Value = "그리드위즈"
$ws.Range("B20").Value = "2024.05.23~05.29"
$ws.Range("C20").Value = "34,000~40,000"
$ws.Range("D20").Value = "-"
$ws.Range("E20").Value = "47600"
$ws.Range("F20").Value = "삼성증권"

